$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 15.67580543370389
$ws.Range("C2").Value = 9.258775060130167
$ws.Range("E2").Value = 11.45936321324002
$ws.Range("F2").Value = 16.86991607391233
$ws.Range("G2").Value = 35.36319745175258
$ws.Range("H2").Value = 15.74866973382392
$ws.Range("I2").Value = 24.05951020855206
$ws.Range("L2").Value = 10.09418631441885
$ws.Range("M2").Value = 15.14087804230605
$ws.Range("N2").Value = 18.3780669924703

$ws.Range("B3").Value = 15.2436934862816
$ws.Range("C3").Value = 8.902507102349899
$ws.Range("E3").Value = 11.47896006517825
$ws.Range("F3").Value = 15.89584955866808
$ws.Range("G3").Value = 35.28745719744976
$ws.Range("H3").Value = 15.79320754334969
$ws.Range("I3").Value = 24.15729846049818
$ws.Range("L3").Value = 10.10330432110615
$ws.Range("M3").Value = 15.06109967948958
$ws.Range("N3").Value = 18.44409461096154

$ws.Range("B4").Value = 14.97524832660323
$ws.Range("C4").Value = 8.674414206122655
$ws.Range("E4").Value = 11.49169694866751
$ws.Range("F4").Value = 15.26997757108491
$ws.Range("G4").Value = 35.25489072364117
$ws.Range("H4").Value = 15.82406261674195
$ws.Range("I4").Value = 24.22308378982684
$ws.Range("L4").Value = 10.11036736379387
$ws.Range("M4").Value = 15.01460747356043
$ws.Range("N4").Value = 18.48658081282228

$ws.Range("B5").Value = 14.86524249034978
$ws.Range("C5").Value = 8.579191148930843
$ws.Range("E5").Value = 11.49706488158269
$ws.Range("F5").Value = 15.00819731993403
$ws.Range("G5").Value = 35.24512630049497
$ws.Range("H5").Value = 15.83751605569234
$ws.Range("I5").Value = 24.25133121434729
$ws.Range("L5").Value = 10.11361419629674
$ws.Range("M5").Value = 14.99630291625177
$ws.Range("N5").Value = 18.50438487234383

$ws.Range("B6").Value = 14.84694439464452
$ws.Range("C6").Value = 8.563244545402945
$ws.Range("E6").Value = 11.49796695922755
$ws.Range("F6").Value = 14.96433081551593
$ws.Range("G6").Value = 35.24371667651607
$ws.Range("H6").Value = 15.83980304710574
$ws.Range("I6").Value = 24.25610847263988
$ws.Range("L6").Value = 10.11417559972676
$ws.Range("M6").Value = 14.9933026110945
$ws.Range("N6").Value = 18.50737089754561

$ws.Range("B7").Value = 14.97376699162597
$ws.Range("C7").Value = 8.673139091618479
$ws.Range("E7").Value = 11.49176862300403
$ws.Range("F7").Value = 15.26647399323137
$ws.Range("G7").Value = 35.25474484017123
$ws.Range("H7").Value = 15.82424049593876
$ws.Range("I7").Value = 24.22345892286654
$ws.Range("L7").Value = 10.11040965896604
$ws.Range("M7").Value = 15.01435799609365
$ws.Range("N7").Value = 18.48681893624113

$ws.Range("B8").Value = 15.5275635575779
$ws.Range("C8").Value = 9.137915564971079
$ws.Range("E8").Value = 11.46597428648395
$ws.Range("F8").Value = 16.53996406344768
$ws.Range("G8").Value = 35.33419184084563
$ws.Range("H8").Value = 15.76329664879837
$ws.Range("I8").Value = 24.09203281085046
$ws.Range("L8").Value = 10.09702644248663
$ws.Range("M8").Value = 15.11286234868971
$ws.Range("N8").Value = 18.40043054442917

$ws.Range("B9").Value = 16.5816128298008
$ws.Range("C9").Value = 9.972353802657922
$ws.Range("E9").Value = 11.42096036954718
$ws.Range("F9").Value = 19.00274580682531
$ws.Range("G9").Value = 35.60027557812943
$ws.Range("H9").Value = 15.67173541820436
$ws.Range("I9").Value = 23.88008490864856
$ws.Range("L9").Value = 10.08238695538111
$ws.Range("M9").Value = 15.32511632525327
$ws.Range("N9").Value = 18.24638632120674

$ws.Range("B10").Value = 17.32775954361523
$ws.Range("C10").Value = 10.53526882486435
$ws.Range("E10").Value = 11.39125575413604
$ws.Range("F10").Value = 20.67494806633232
$ws.Range("G10").Value = 35.86226152184788
$ws.Range("H10").Value = 15.62164821052069
$ws.Range("I10").Value = 23.75256036138153
$ws.Range("L10").Value = 10.07868226031747
$ws.Range("M10").Value = 15.49172641037755
$ws.Range("N10").Value = 18.14247785819829

$ws.Range("B11").Value = 17.65943264698318
$ws.Range("C11").Value = 10.7799377760987
$ws.Range("E11").Value = 11.37846763178234
$ws.Range("F11").Value = 21.3917225636224
$ws.Range("G11").Value = 35.9956236961842
$ws.Range("H11").Value = 15.60262164421165
$ws.Range("I11").Value = 23.70072848125676
$ws.Range("L11").Value = 10.07852135412219
$ws.Range("M11").Value = 15.56961761611875
$ws.Range("N11").Value = 18.09719884080203

$ws.Range("B12").Value = 17.78379339549703
$ws.Range("C12").Value = 10.87091207741732
$ws.Range("E12").Value = 11.37372886722042
$ws.Range("F12").Value = 21.65686569030329
$ws.Range("G12").Value = 36.04813471910855
$ws.Range("H12").Value = 15.59595931384324
$ws.Range("I12").Value = 23.68199464907569
$ws.Range("L12").Value = 10.07867890708835
$ws.Range("M12").Value = 15.59939545329854
$ws.Range("N12").Value = 18.08033740003316

$ws.Range("B13").Value = 17.75706690389656
$ws.Range("C13").Value = 10.85139422154014
$ws.Range("E13").Value = 11.3747448336238
$ws.Range("F13").Value = 21.60004134736742
$ws.Range("G13").Value = 36.03673667490856
$ws.Range("H13").Value = 15.59736999993279
$ws.Range("I13").Value = 23.68598948297033
$ws.Range("L13").Value = 10.07863527044126
$ws.Range("M13").Value = 15.59297003498579
$ws.Range("N13").Value = 18.08395617026673

$ws.Range("B14").Value = 17.6696893227902
$ws.Range("C14").Value = 10.78745607650901
$ws.Range("E14").Value = 11.37807569225335
$ws.Range("F14").Value = 21.4136618050453
$ws.Range("G14").Value = 35.99990372046445
$ws.Range("H14").Value = 15.60206264354629
$ws.Range("I14").Value = 23.69916929580667
$ws.Range("L14").Value = 10.07852994176632
$ws.Range("M14").Value = 15.57206191724462
$ws.Range("N14").Value = 18.09580594092661

$ws.Range("B15").Value = 17.61600355918748
$ws.Range("C15").Value = 10.74807279282374
$ws.Range("E15").Value = 11.38012944840548
$ws.Range("F15").Value = 21.29868154950795
$ws.Range("G15").Value = 35.9776032244086
$ws.Range("H15").Value = 15.60500774981031
$ws.Range("I15").Value = 23.70735885915681
$ws.Range("L15").Value = 10.0784938552527
$ws.Range("M15").Value = 15.55929123848875
$ws.Range("N15").Value = 18.10310130968307

$ws.Range("B16").Value = 17.30591724883741
$ws.Range("C16").Value = 10.51904639349711
$ws.Range("E16").Value = 11.39210603694073
$ws.Range("F16").Value = 20.62722412089977
$ws.Range("G16").Value = 35.85382894618274
$ws.Range("H16").Value = 15.62296745815799
$ws.Range("I16").Value = 23.75607248913172
$ws.Range("L16").Value = 10.07872339079724
$ws.Range("M16").Value = 15.48667660358409
$ws.Range("N16").Value = 18.14547686645155

$ws.Range("B17").Value = 17.113613391312
$ws.Range("C17").Value = 10.37559737796404
$ws.Range("E17").Value = 11.39963861263492
$ws.Range("F17").Value = 20.20408069597325
$ws.Range("G17").Value = 35.78151034842868
$ws.Range("H17").Value = 15.63494935566708
$ws.Range("I17").Value = 23.78754317155828
$ws.Range("L17").Value = 10.07925414889599
$ws.Range("M17").Value = 15.44265367769028
$ws.Range("N17").Value = 18.17198148703697

$ws.Range("B18").Value = 17.00228212344755
$ws.Range("C18").Value = 10.2920177096257
$ws.Range("E18").Value = 11.40403938176419
$ws.Range("F18").Value = 19.95656407809801
$ws.Range("G18").Value = 35.74125170608153
$ws.Range("H18").Value = 15.64219466008303
$ws.Range("I18").Value = 23.80622563263935
$ws.Range("L18").Value = 10.07970295705026
$ws.Range("M18").Value = 15.41753167948855
$ws.Range("N18").Value = 18.18741358056592

$ws.Range("B19").Value = 16.96446696059058
$ws.Range("C19").Value = 10.26353623492015
$ws.Range("E19").Value = 11.4055411370775
$ws.Range("F19").Value = 19.87204792380568
$ws.Range("G19").Value = 35.7278513059011
$ws.Range("H19").Value = 15.64470846401354
$ws.Range("I19").Value = 23.81265088360725
$ws.Range("L19").Value = 10.07987958772284
$ws.Range("M19").Value = 15.40906054561427
$ws.Range("N19").Value = 18.19267083874923

$ws.Range("B20").Value = 17.13416022942429
$ws.Range("C20").Value = 10.39097895955796
$ws.Range("E20").Value = 11.39882969826961
$ws.Range("F20").Value = 20.24955283636154
$ws.Range("G20").Value = 35.78907060353078
$ws.Range("H20").Value = 15.63363724619602
$ws.Range("I20").Value = 23.78413285786697
$ws.Range("L20").Value = 10.07918279889552
$ws.Range("M20").Value = 15.44731955746792
$ws.Range("N20").Value = 18.16914064602691

$ws.Range("B21").Value = 17.69538871393686
$ws.Range("C21").Value = 10.80628202404964
$ws.Range("E21").Value = 11.37709452373719
$ws.Range("F21").Value = 21.46857628470577
$ws.Range("G21").Value = 36.01066817122238
$ws.Range("H21").Value = 15.60066955736025
$ws.Range("I21").Value = 23.69527376783762
$ws.Range("L21").Value = 10.0785549558983
$ws.Range("M21").Value = 15.57819564111009
$ws.Range("N21").Value = 18.09231765887622

$ws.Range("B22").Value = 18.05492714416807
$ws.Range("C22").Value = 11.06792138104012
$ws.Range("E22").Value = 11.36349429681842
$ws.Range("F22").Value = 22.22866616901552
$ws.Range("G22").Value = 36.1671918091186
$ws.Range("H22").Value = 15.58228695164367
$ws.Range("I22").Value = 23.64241115746134
$ws.Range("L22").Value = 10.07941772547293
$ws.Range("M22").Value = 15.66536572244525
$ws.Range("N22").Value = 18.04376847042213

$ws.Range("B23").Value = 17.86373607049641
$ws.Range("C23").Value = 10.92918556073526
$ws.Range("E23").Value = 11.37069776539591
$ws.Range("F23").Value = 21.82633154458858
$ws.Range("G23").Value = 36.08259298042076
$ws.Range("H23").Value = 15.59180794049896
$ws.Range("I23").Value = 23.67014635458972
$ws.Range("L23").Value = 10.07884102009462
$ws.Range("M23").Value = 15.6186984483056
$ws.Range("N23").Value = 18.06952871952959

$ws.Range("B24").Value = 17.12487340290657
$ws.Range("C24").Value = 10.38402839783409
$ws.Range("E24").Value = 11.39919518980155
$ws.Range("F24").Value = 20.22900810905287
$ws.Range("G24").Value = 35.78564850354192
$ws.Range("H24").Value = 15.63422933973828
$ws.Range("I24").Value = 23.78567282460022
$ws.Range("L24").Value = 10.07921460866321
$ws.Range("M24").Value = 15.44520952806187
$ws.Range("N24").Value = 18.1704243850788

$ws.Range("B25").Value = 16.30083581106041
$ws.Range("C25").Value = 9.755220612857329
$ws.Range("E25").Value = 11.43254454117894
$ws.Range("F25").Value = 18.34778573295695
$ws.Range("G25").Value = 35.51654266472634
$ws.Range("H25").Value = 15.69349758815241
$ws.Range("I25").Value = 23.93249183398241
$ws.Range("L25").Value = 10.08510706854142
$ws.Range("M25").Value = 15.26574978898244
$ws.Range("N25").Value = 18.28642482364857
